$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 133: clone formatting + values from row 132 (identical B:H values),
# then fix the date in column A.
$ws.Range("A132:H132").Copy()
$ws.Range("A133").PasteSpecial(-4122)
$ws.Range("A132:H132").Copy()
$ws.Range("A133").PasteSpecial(-4163)
$ws.Range("A133").Value = 45504.2916666667

# Row 134: clone formatting + values from row 132 again, then fix the
# date in column A and the volume in column B.
$ws.Range("A132:H132").Copy()
$ws.Range("A134").PasteSpecial(-4122)
$ws.Range("A132:H132").Copy()
$ws.Range("A134").PasteSpecial(-4163)
$ws.Range("A134").Value = 45505.609212963
$ws.Range("B134").Value = 900
